# Insert a new "creator" row above row 13 (this naturally shifts the
# existing rows 13-90 down to 14-91, matching the diff), then populate
# the newly inserted row 13 with the new creator's details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 - pushes old row13..row90 -> row14..row91
$ws.Rows.Item(13).Insert()

# Fill the new row 13 with the new contributor entry
$ws.Range("A13").Value = "dct:creator"
$ws.Range("B13").Value = "https://orcid.org/0000-0001-5725-8297"
$ws.Range("C13").Value = "Hannah Random Found"
